$d = $word.ActiveDocument

$replacements = @(
    @("306×2=", "595×6="),
    @("285×2=", "812×9="),
    @("808×7=", "656×6="),
    @("581×8=", "960×9="),
    @("711×7=", "826×4="),
    @("728×5=", "556×8="),
    @("225×9=", "753×7="),
    @("743×4=", "539×5="),
    @("952×5=", "730×5="),
    @("196×2=", "178×7="),
    @("115×5=", "338×9="),
    @("564×7=", "172×4="),
    @("290×6=", "954×6="),
    @("573×8=", "846×9="),
    @("688×3=", "634×3="),
    @("261×9=", "299×9="),
    @("213×7=", "634×3="),
    @("171×7=", "977×9="),
    @("578×8=", "602×3="),
    @("147×7=", "498×3="),
    @("847×8=", "937×3="),
    @("706×4=", "663×3="),
    @("268×2=", "243×2="),
    @("521×6=", "361×5="),
    @("394×4=", "794×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
